$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2; this shifts the Nifemi/Test User/ace rows down
# to rows 3/4/5 respectively (matches dimension growing from A1:E4 to A1:E5).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the "Mimi" record.
$ws.Cells.Item(2, 1).Value = "Mimi"

# Force the phone number to be stored as text so the leading zero is kept.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "08035987700"
$ws.Cells.Item(2, 2).ClearFormats()

$ws.Cells.Item(2, 3).Value = "smilingmiriam@gmail.com"
$ws.Cells.Item(2, 4).Value = 12
$ws.Cells.Item(2, 5).Value = "2025-09-18T14:15:27.058764+00:00"

# The "Nifemi  Spectro" record (now on row 3) gets an updated id and created_at.
$ws.Cells.Item(3, 4).Value = 11
$ws.Cells.Item(3, 5).Value = "2025-09-18T09:47:44.947181+00:00"
